$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking Price cells to remain as literal text
# (keeps trailing zeros / exact decimal string instead of Excel coercing to a Double)
$textCells = @("D4","D5","D6","D10","D11","D12","D16","D19","D21","D22","D25","D27","D28","D30","D31","D33","D34","D35","D36","D39","D41","D42","D43","D45","D47","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '42.922.62'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").Value = '2.573.17'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '302.33'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '97.49'
$ws.Range("E6").Value = '  +1.85%  '
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("D10").Value = '36.38'
$ws.Range("E10").Value = '  -1.46%  '
$ws.Range("D11").Value = '0.0810'
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = '7.73'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("E13").Value = '  +6.01%  '
$ws.Range("D14").Value = '2.560.66'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '14.40'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").Value = '42.992.11'
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("E18").Value = '  +2.04%  '
$ws.Range("D19").Value = '12.89'
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '72.06'
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("D22").Value = '254.70'
$ws.Range("E22").Value = '  -3.52%  '
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("E24").Value = '  -2.71%  '
$ws.Range("D25").Value = '28.82'
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '10.25'
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").Value = '37.78'
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("D30").Value = '6.03'
$ws.Range("E30").Value = '  -1.70%  '
$ws.Range("D31").Value = '155.33'
$ws.Range("E31").Value = '  +2.53%  '
$ws.Range("E32").Value = '  -3.58%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '2.17'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '0.0806'
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").Value = '18.31'
$ws.Range("E36").Value = '  +11.94%  '
$ws.Range("E37").Value = '  -2.52%  '
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").Value = '23.10'
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("E40").Value = '  -3.67%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0311'
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value = '2.07'
$ws.Range("E42").Value = '  +26.47%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '3.88'
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = '2.075.23'
$ws.Range("E44").Value = '  +2.38%  '
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").Value = '85.42'
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("D48").Value = '76.39'
$ws.Range("E48").Value = '  +10.80%  '
$ws.Range("D49").Value = '106.72'
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("D50").Value = '2.822.93'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("E51").Value = '  +1.78%  '
